$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zipcode")

# Update the B4 answer value (Columbus, OH zip) to reflect new submitted value
$ws.Range("B4").Value = 432215

# Reveal the hidden hint cell for the failed/incorrect answer
$ws.Range("D4").Value = "// wrong"

# Update the active selection to match the cell being worked on
$ws.Range("F6").Select()
